$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "tipoScheda" tipologica (row 81: A81 = "tipoScheda") was removed from
# the sheet. Deleting the entire row shifts every row below it up by one,
# which is exactly what the target workbook shows (old row 82
# "tipoServizioSoggetto" becomes the new row 81, old row 83 "tipoSoggetto"
# becomes the new row 82, and the sheet's used range shrinks from A1:B83 to
# A1:B82). Excel also prunes the now-unreferenced "tipoScheda" shared
# string when it rewrites sharedStrings.xml on save.
$ws.Rows(81).Delete()

# Match the final view/selection state recorded in the saved file: the
# window was scrolled down and the (now renumbered) row 81 was selected in
# full before the file was saved.
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows("81:81").Select()
